$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.028770432315545
$ws.Range("D2").Value = 1.030023074229415
$ws.Range("E2").Value = 1.038020203507392
$ws.Range("F2").Value = 1.047383162933067
$ws.Range("I2").Value = 1.034297704200607
$ws.Range("J2").Value = 1.033920907391606
$ws.Range("K2").Value = 1.032835288492301
$ws.Range("L2").Value = 1.040809427505572
$ws.Range("M2").Value = 1.050145954202725
$ws.Range("N2").Value = 1.015327804369364

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.029625315280833
$ws.Range("D3").Value = 1.030751398150041
$ws.Range("E3").Value = 1.038788247179736
$ws.Range("F3").Value = 1.048256171589697
$ws.Range("I3").Value = 1.034422360581718
$ws.Range("J3").Value = 1.034417024629202
$ws.Range("K3").Value = 1.033372095355162
$ws.Range("L3").Value = 1.041387509376476
$ws.Range("M3").Value = 1.050830637291421
$ws.Range("N3").Value = 1.015491903973577

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.030178742735585
$ws.Range("D4").Value = 1.031223209421887
$ws.Range("E4").Value = 1.039285825239878
$ws.Range("F4").Value = 1.048821652621188
$ws.Range("I4").Value = 1.034501114938722
$ws.Range("J4").Value = 1.034737663855557
$ws.Range("K4").Value = 1.033719314458938
$ws.Range("L4").Value = 1.041761508724116
$ws.Range("M4").Value = 1.051273619319814
$ws.Range("N4").Value = 1.01559794561833

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.030411464611268
$ws.Range("D5").Value = 1.031421685698871
$ws.Range("E5").Value = 1.039495149521794
$ws.Range("F5").Value = 1.049059519283114
$ws.Range("I5").Value = 1.034533765913239
$ws.Range("J5").Value = 1.034872367902098
$ws.Range("K5").Value = 1.033865252689696
$ws.Range("L5").Value = 1.041918722681169
$ws.Range("M5").Value = 1.051459834470228
$ws.Range("N5").Value = 1.015642491137692

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.030450543169716
$ws.Range("D6").Value = 1.031455018125121
$ws.Range("E6").Value = 1.039530304292716
$ws.Range("F6").Value = 1.049099466205919
$ws.Range("I6").Value = 1.034539221309621
$ws.Range("J6").Value = 1.034894979831871
$ws.Range("K6").Value = 1.033889754404787
$ws.Range("L6").Value = 1.041945118653815
$ws.Range("M6").Value = 1.051491099915762
$ws.Range("N6").Value = 1.015649968496822

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.030181852138994
$ws.Range("D7").Value = 1.0312258609755
$ws.Range("E7").Value = 1.03928862168361
$ws.Range("F7").Value = 1.048824830465562
$ws.Range("I7").Value = 1.03450155302088
$ws.Range("J7").Value = 1.034739464143422
$ws.Range("K7").Value = 1.033721264622849
$ws.Range("L7").Value = 1.041763609486882
$ws.Range("M7").Value = 1.051276107594785
$ws.Range("N7").Value = 1.015598540973427

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.029059289063148
$ws.Range("D8").Value = 1.030269102476587
$ws.Range("E8").Value = 1.038279641664567
$ws.Range("F8").Value = 1.047678077829217
$ws.Range("I8").Value = 1.034340226655098
$ws.Range("J8").Value = 1.034088650821
$ws.Range("K8").Value = 1.033016731158947
$ws.Range("L8").Value = 1.041004804196993
$ws.Range("M8").Value = 1.050377356198563
$ws.Range("N8").Value = 1.015383291685266

# Row 9
$ws.Range("B9").Value = 1.019999999999999
$ws.Range("C9").Value = 1.027083255166386
$ws.Range("D9").Value = 1.028587361070658
$ws.Range("E9").Value = 1.036506381558111
$ws.Range("F9").Value = 1.04566192382805
$ws.Range("I9").Value = 1.034041385710964
$ws.Range("J9").Value = 1.03293896676841
$ws.Range("K9").Value = 1.031774310966391
$ws.Range("L9").Value = 1.039667308338445
$ws.Range("M9").Value = 1.048793301989469
$ws.Range("N9").Value = 1.015002928817676

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.025767376920992
$ws.Range("D10").Value = 1.027469119555324
$ws.Range("E10").Value = 1.035327466590472
$ws.Range("F10").Value = 1.044321004295196
$ws.Range("I10").Value = 1.033832419513845
$ws.Range("J10").Value = 1.03217065790154
$ws.Range("K10").Value = 1.030945477442276
$ws.Range("L10").Value = 1.038775469769609
$ws.Range("M10").Value = 1.04773713398384
$ws.Range("N10").Value = 1.014748663129978

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.025197956028747
$ws.Range("D11").Value = 1.026985622619498
$ws.Range("E11").Value = 1.03481777805043
$ws.Range("F11").Value = 1.04374114933271
$ws.Range("I11").Value = 1.033739635819325
$ws.Range("J11").Value = 1.031837547918842
$ws.Range("K11").Value = 1.030586468849701
$ws.Range("L11").Value = 1.038389269839488
$ws.Range("M11").Value = 1.047279790411675
$ws.Range("N11").Value = 1.014638405018186

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.024986503753157
$ws.Range("D12").Value = 1.026806138580028
$ws.Range("E12").Value = 1.034628577381841
$ws.Range("F12").Value = 1.043525883208603
$ws.Range("I12").Value = 1.03370482712366
$ws.Range("J12").Value = 1.031713753104108
$ws.Range("K12").Value = 1.030453100644707
$ws.Range("L12").Value = 1.038245815178964
$ws.Range("M12").Value = 1.047109911861579
$ws.Range("N12").Value = 1.014597426735598

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.025031858426812
$ws.Range("D13").Value = 1.026844633589538
$ws.Range("E13").Value = 1.034669156093464
$ws.Range("F13").Value = 1.043572053161618
$ws.Range("I13").Value = 1.033712309292001
$ws.Range("J13").Value = 1.031740310340161
$ws.Range("K13").Value = 1.030481709299539
$ws.Range("L13").Value = 1.03827658679559
$ws.Range("M13").Value = 1.047146351397126
$ws.Range("N13").Value = 1.014606217773239

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.02518047616698
$ws.Range("D14").Value = 1.026970784194564
$ws.Range("E14").Value = 1.034802136198974
$ws.Range("F14").Value = 1.043723352953506
$ws.Range("I14").Value = 1.033736765548618
$ws.Range("J14").Value = 1.031827316277535
$ws.Range("K14").Value = 1.030575444919741
$ws.Range("L14").Value = 1.038377411881835
$ws.Range("M14").Value = 1.047265748204507
$ws.Range("N14").Value = 1.014635018217659

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.025272051925219
$ws.Range("D15").Value = 1.027048524141599
$ws.Range("E15").Value = 1.034884085630468
$ws.Range("F15").Value = 1.043816589416861
$ws.Range("I15").Value = 1.033751788204422
$ws.Range("J15").Value = 1.031880915204637
$ws.Range("K15").Value = 1.03063319639174
$ws.Range("L15").Value = 1.038439533214059
$ws.Range("M15").Value = 1.047339312461921
$ws.Range("N15").Value = 1.014652760019628

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.02580517528357
$ws.Range("D16").Value = 1.027501222740467
$ws.Range("E16").Value = 1.03536130972174
$ws.Range("F16").Value = 1.044359503782281
$ws.Range("I16").Value = 1.033838528876671
$ws.Range("J16").Value = 1.032192756420138
$ws.Range("K16").Value = 1.030969301283388
$ws.Range("L16").Value = 1.038801100116833
$ws.Range("M16").Value = 1.047767486172961
$ws.Range("N16").Value = 1.014755977276807

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.026139687710057
$ws.Range("D17").Value = 1.027785379657907
$ws.Range("E17").Value = 1.035660872371903
$ws.Range("F17").Value = 1.044700267616107
$ws.Range("I17").Value = 1.0338923238419
$ws.Range("J17").Value = 1.032388252688119
$ws.Range("K17").Value = 1.031180100295521
$ws.Range("L17").Value = 1.039027894981916
$ws.Range("M17").Value = 1.048036065114428
$ws.Range("N17").Value = 1.014820680370147

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.026334838082958
$ws.Range("D18").Value = 1.027951191874872
$ws.Range("E18").Value = 1.03583567824832
$ws.Range("F18").Value = 1.044899103763332
$ws.Range("I18").Value = 1.033923479566808
$ws.Range("J18").Value = 1.032502241054618
$ws.Range("K18").Value = 1.031303044298647
$ws.Range("L18").Value = 1.039160177800643
$ws.Range("M18").Value = 1.048192721030252
$ws.Range("N18").Value = 1.014858405181726

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.026401385198254
$ws.Range("D19").Value = 1.028007741087163
$ws.Range("E19").Value = 1.035895295381066
$ws.Range("F19").Value = 1.044966914309346
$ws.Range("I19").Value = 1.033934065188532
$ws.Range("J19").Value = 1.032541101098955
$ws.Range("K19").Value = 1.031344963036384
$ws.Range("L19").Value = 1.039205282307984
$ws.Range("M19").Value = 1.048246136322876
$ws.Range("N19").Value = 1.014871265740807

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.026103794068596
$ws.Range("D20").Value = 1.027754885233374
$ws.Range("E20").Value = 1.035628724255645
$ws.Range("F20").Value = 1.044663699176154
$ws.Range("I20").Value = 1.03388657510435
$ws.Range("J20").Value = 1.032367282051485
$ws.Range("K20").Value = 1.031157484735773
$ws.Range("L20").Value = 1.039003562317697
$ws.Range("M20").Value = 1.048007249293543
$ws.Range("N20").Value = 1.014813739931149

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.025136710389267
$ws.Range("D21").Value = 1.026933632983522
$ws.Range("E21").Value = 1.034762973536114
$ws.Range("F21").Value = 1.043678795668521
$ws.Range("I21").Value = 1.033729573299284
$ws.Range("J21").Value = 1.031801696924914
$ws.Range("K21").Value = 1.030547842558746
$ws.Range("L21").Value = 1.038347721473412
$ws.Range("M21").Value = 1.04723058882946
$ws.Range("N21").Value = 1.014626537850124

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.024528991195917
$ws.Range("D22").Value = 1.026417906251883
$ws.Range("E22").Value = 1.03421933886074
$ws.Range("F22").Value = 1.043060230358385
$ws.Range("I22").Value = 1.033628866116159
$ws.Range("J22").Value = 1.031445727594539
$ws.Range("K22").Value = 1.030164441928205
$ws.Range("L22").Value = 1.03793535272457
$ws.Range("M22").Value = 1.04674226782036
$ws.Range("N22").Value = 1.014508700661159

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.02485112337735
$ws.Range("D23").Value = 1.026691242727906
$ws.Range("E23").Value = 1.034507463222978
$ws.Range("F23").Value = 1.043388078217118
$ws.Range("I23").Value = 1.033682441620015
$ws.Range("J23").Value = 1.031634467714916
$ws.Range("K23").Value = 1.030367698361434
$ws.Range("L23").Value = 1.038153958215369
$ws.Range("M23").Value = 1.047001135849863
$ws.Range("N23").Value = 1.014571181112692

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.02612001275235
$ws.Range("D24").Value = 1.027768664140387
$ws.Range("E24").Value = 1.0356432503708
$ws.Range("F24").Value = 1.044680222650515
$ws.Range("I24").Value = 1.033889173397324
$ws.Range("J24").Value = 1.03237675790725
$ws.Range("K24").Value = 1.03116770376941
$ws.Range("L24").Value = 1.03901455721004
$ws.Range("M24").Value = 1.048020269927867
$ws.Range("N24").Value = 1.01481687606444

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.02759385333706
$ws.Range("D25").Value = 1.029021624467647
$ws.Range("E25").Value = 1.036964245076014
$ws.Range("F25").Value = 1.04618259495331
$ws.Range("I25").Value = 1.034120363637971
$ws.Range("J25").Value = 1.033236519724684
$ws.Range("K25").Value = 1.032095609751269
$ws.Range("L25").Value = 1.040013119957886
$ws.Range("M25").Value = 1.049202848246224
$ws.Range("N25").Value = 1.015101385354297

